# Add a new dialogue triplet (English / Russian / "corrupted" Russian)
# to the Кроганк (Croagunk) character sheet.
#
# Before the edit, row 38 (B=234) was the last row of a "line group" and
# had no bottom border (style 4/5, same as rows 36-37). The edit:
#   1. Turns row 38 into the closing (bordered) row of its group
#      (style 6 on A/B, style 7 on C/D/E - matching the border style
#      already used by every other group's closing row, e.g. row 35).
#   2. Appends a brand-new group of three rows (39-41, B=198/201/204)
#      holding the new dialogue, using the normal (non-bordered) row
#      style already used by rows 36-37 etc.
#   3. Moves the active selection to D35 (matches the author's final
#      cursor position in the diff).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Give row 38 the "closing row of group" bordered format --------
# Row 35 is an existing closing row (A35/B35 style 6, C35:E35 style 7);
# copy its formatting onto row 38 so the exact same style indices are
# reused instead of minting new ones.
$ws.Range("A35:E35").Copy()
$ws.Range("A38:E38").PasteSpecial(-4122)   # -4122 = xlPasteFormats

# --- 2. Add the new group: rows 39, 40, 41 -----------------------------
# Rows 36/37 are existing "plain" (non-closing) rows - copy their format
# (style 4 on B, style 5 on C/D/E, no border) onto the three new rows.
$ws.Range("B36:E36").Copy()
$ws.Range("B39:E41").PasteSpecial(-4122)   # -4122 = xlPasteFormats

# Row numbers (column B) for the new group.
$ws.Range("B39").Value2 = 198
$ws.Range("B40").Value2 = 201
$ws.Range("B41").Value2 = 204

# Fill column-by-column (C, then D, then E) so brand-new shared strings
# get appended to xl/sharedStrings.xml in the same order the author's
# saved file has them.

# Column C (English). C39 reuses the shared string already used by C36.
$ws.Range("C39").Value2 = $ws.Range("C36").Value2
$ws.Range("C40").Value2 = " The graduation exam will put\nyou through the wringer."
$ws.Range("C41").Value2 = " Don\'t be flattened! Meh-heh-heh."

# Column D (Russian translation).
$ws.Range("D39").Value2 = " Хе-хе-хе. Эй, вы двое..."
$ws.Range("D40").Value2 = " Выпускной экзамен проверит вас\nна прочность."
$ws.Range("D41").Value2 = " Смотрите, не прогнитесь!\nХе-хе-хе!"

# Column E ("converted"/cipher string).
$ws.Range("E39").Value2 = " Öå-öå-öå. Üê, âú äâïå..."
$ws.Range("E40").Value2 = " Âúðôòëîïê üëèàíåî ðñïâåñéó âàò\nîà ðñïœîïòóû."
$ws.Range("E41").Value2 = " Òíïóñéóå, îå ðñïãîéóåòû!\nÖå-öå-öå!"

# Row 40 keeps the taller (wrapped, two-line) row height, same as other
# two-line rows in the sheet.
$ws.Rows.Item(40).RowHeight = 21.6

# --- 3. Move the selection to match the author's final cursor ---------
$ws.Range("D35").Select()
